$wb = $excel.ActiveWorkbook

# ---- Sheet1: LP1912 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("A2").Value = 'Última actualización: 07:51:34'
$ws.Range("A3").Value = 'Total filas: 107'
$rows1 = @(
    @('01:10:32','01:12','215_ALUAR',2,'LP1912'),
    @('01:55:51','01:58','14_ABASTO',3,'LP1912'),
    @('02:48:52','02:57','215_ALUAR',9,'LP1912'),
    @('02:58:51','02:58','215_ALUAR',0,'LP1912'),
    @('01:55:51','03:12','215_ALUAR',77,'LP1912'),
    @('03:35:49','03:38','14_ABASTO',3,'LP1912'),
    @('02:58:51','03:48','14_ABASTO',50,'LP1912'),
    @('02:21:47','03:56','14_ABASTO',95,'LP1912'),
    @('04:01:13','04:01','81_EL PELIGRO',0,'LP1912'),
    @('03:35:49','04:45','215A_EL PATO',70,'LP1912'),
    @('04:35:25','04:46','215A_EL PATO',11,'LP1912'),
    @('04:48:57','04:53','11_ETCHEVERRY',5,'LP1912'),
    @('04:48:57','05:14','14_ABASTO',26,'LP1912'),
    @('04:35:25','05:15','14_ABASTO',40,'LP1912'),
    @('04:56:11','05:16','17_ROMERO',20,'LP1912'),
    @('04:48:57','05:21','23_HERNANDEZ',33,'LP1912'),
    @('04:56:11','05:22','23_HERNANDEZ',26,'LP1912'),
    @('05:21:16','05:26','23_HERNANDEZ',5,'LP1912'),
    @('03:35:49','05:28','14_ABASTO',113,'LP1912'),
    @('04:48:57','05:34','215B_EL PATO',46,'LP1912'),
    @('05:21:16','05:35','215B_EL PATO',14,'LP1912'),
    @('04:01:13','05:37','14_ABASTO',96,'LP1912'),
    @('05:21:16','05:46','15_ABASTO',25,'LP1912'),
    @('05:52:07','05:54','10_OLMOS',2,'LP1912'),
    @('05:21:16','06:04','16_SANTA ANA',43,'LP1912'),
    @('05:21:16','06:11','215A_EL PATO',50,'LP1912'),
    @('05:52:07','06:12','215A_EL PATO',20,'LP1912'),
    @('04:48:57','06:13','225_HARAS DEL SUR',85,'LP1912'),
    @('05:52:07','06:14','225_HARAS DEL SUR',22,'LP1912'),
    @('04:48:57','06:20','26_HERNANDEZ',92,'LP1912'),
    @('06:21:22','06:21','26_HERNANDEZ',0,'LP1912'),
    @('04:48:57','06:26','23_HERNANDEZ',98,'LP1912'),
    @('05:52:07','06:27','23_HERNANDEZ',35,'LP1912'),
    @('06:21:22','06:29','23_HERNANDEZ',8,'LP1912'),
    @('06:21:22','06:29','86_EST CHICA-ESC AGRARIA',8,'LP1912'),
    @('05:52:07','06:30','86_EST CHICA-ESC AGRARIA',38,'LP1912'),
    @('06:21:22','06:31','16_SANTA ANA',10,'LP1912'),
    @('04:48:57','06:43','225_C ROCA-H SUR',115,'LP1912'),
    @('06:21:22','06:44','225_C ROCA-H SUR',23,'LP1912'),
    @('06:21:22','06:46','215C_EL PATO',25,'LP1912'),
    @('05:52:07','06:47','215C_EL PATO',55,'LP1912'),
    @('06:59:37','06:59','14_ABASTO',0,'LP1912'),
    @('06:59:37','06:59','23_HERNANDEZ',0,'LP1912'),
    @('05:52:07','07:00','14_ABASTO',68,'LP1912'),
    @('06:49:33','07:01','16_SANTA ANA',12,'LP1912'),
    @('06:49:33','07:04','23_HERNANDEZ',15,'LP1912'),
    @('05:52:07','07:05','23_HERNANDEZ',73,'LP1912'),
    @('06:59:37','07:05','15_ABASTO',6,'LP1912'),
    @('06:59:37','07:07','225_GOMEZ',8,'LP1912'),
    @('06:59:37','07:11','215A_EL PATO',12,'LP1912'),
    @('05:52:07','07:12','215A_EL PATO',80,'LP1912'),
    @('06:59:37','07:15','11_ETCHEVERRY',16,'LP1912'),
    @('05:52:07','07:16','11_ETCHEVERRY',84,'LP1912'),
    @('06:59:37','07:16','16_SANTA ANA',17,'LP1912'),
    @('06:59:37','07:21','26_HERNANDEZ',22,'LP1912'),
    @('06:59:37','07:23','10_OLMOS',24,'LP1912'),
    @('07:28:14','07:30','11_ETCHEVERRY',2,'LP1912'),
    @('06:59:37','07:31','11_ETCHEVERRY',32,'LP1912'),
    @('06:59:37','07:31','16_SANTA ANA',32,'LP1912'),
    @('05:52:07','07:32','11_ETCHEVERRY',100,'LP1912'),
    @('07:28:14','07:32','84_COLONIA URQUIZA-ESC 49',4,'LP1912'),
    @('05:52:07','07:32','16_SANTA ANA',100,'LP1912'),
    @('07:28:14','07:35','23_HERNANDEZ',7,'LP1912'),
    @('06:59:37','07:36','27_EL RETIRO',37,'LP1912'),
    @('07:28:14','07:37','27_EL RETIRO',9,'LP1912'),
    @('07:28:14','07:39','10_OLMOS',11,'LP1912'),
    @('06:59:37','07:47','14_ABASTO',48,'LP1912'),
    @('07:28:14','07:47','16_SANTA ANA',19,'LP1912'),
    @('07:28:14','07:48','14_ABASTO',20,'LP1912'),
    @('07:51:34','07:51','215D_EL PATO',0,'LP1912'),
    @('07:51:34','07:51','10_OLMOS',0,'LP1912'),
    @('07:28:14','07:55','10_OLMOS',27,'LP1912'),
    @('07:28:14','08:00','23_HERNANDEZ',32,'LP1912'),
    @('07:51:34','08:01','23_HERNANDEZ',10,'LP1912'),
    @('07:51:34','08:03','11_ETCHEVERRY',12,'LP1912'),
    @('06:59:37','08:06','23_HERNANDEZ',67,'LP1912'),
    @('07:51:34','08:10','16_SANTA ANA',19,'LP1912'),
    @('07:28:14','08:11','16_SANTA ANA',43,'LP1912'),
    @('07:51:34','08:12','15_ABASTO',21,'LP1912'),
    @('07:51:34','08:13','10_OLMOS',22,'LP1912'),
    @('07:51:34','08:21','26_HERNANDEZ',30,'LP1912'),
    @('07:51:34','08:22','16_P MOR-SANTA ANA',31,'LP1912'),
    @('07:51:34','08:23','215B_EL PATO',32,'LP1912'),
    @('07:28:14','08:23','16_P MOR-SANTA ANA',55,'LP1912'),
    @('07:51:34','08:27','84_COLONIA URQUIZA-ESC 49',36,'LP1912'),
    @('07:51:34','08:30','23_HERNANDEZ',39,'LP1912'),
    @('07:51:34','08:42','81_EL PELIGRO',51,'LP1912'),
    @('07:51:34','08:43','14_ABASTO',52,'LP1912'),
    @('07:28:14','08:44','14_ABASTO',76,'LP1912'),
    @('07:51:34','08:54','17_ROMERO',63,'LP1912'),
    @('07:51:34','09:01','215A_EL PATO',70,'LP1912'),
    @('07:28:14','09:02','215A_EL PATO',94,'LP1912'),
    @('07:51:34','09:03','11_ETCHEVERRY',72,'LP1912'),
    @('07:51:34','09:10','16_P MOR-SANTA ANA',79,'LP1912'),
    @('07:28:14','09:11','16_P MOR-SANTA ANA',103,'LP1912'),
    @('07:51:34','09:16','27_EL RETIRO',85,'LP1912'),
    @('07:28:14','09:17','27_EL RETIRO',109,'LP1912'),
    @('07:51:34','09:21','26_HERNANDEZ',90,'LP1912'),
    @('07:51:34','09:22','17_ROMERO',91,'LP1912'),
    @('07:28:14','09:23','17_ROMERO',115,'LP1912'),
    @('07:51:34','09:23','11_ETCHEVERRY',92,'LP1912'),
    @('07:28:14','09:24','11_ETCHEVERRY',116,'LP1912'),
    @('07:51:34','09:32','15_ABASTO',101,'LP1912'),
    @('07:51:34','09:33','10_OLMOS',102,'LP1912'),
    @('07:51:34','09:34','16_SANTA ANA',103,'LP1912'),
    @('07:51:34','09:42','215C_EL PATO',111,'LP1912'),
    @('07:51:34','09:43','14_ABASTO',112,'LP1912')
)
$n1 = $rows1.Count
$arr1 = New-Object 'object[,]' $n1,5
for ($i = 0; $i -lt $n1; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $arr1[$i,$j] = $rows1[$i][$j]
    }
}
$ws.Range($ws.Cells.Item(6,1), $ws.Cells.Item(5+$n1,5)).Value = $arr1

# ---- Sheet2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A2").Value = 'Última actualización: 07:51:34'
$ws.Range("A3").Value = 'Total filas: 19'
$rows2 = @(
    @('01:10:32','01:12','215_ALUAR',2,'LP1912'),
    @('02:48:52','02:57','215_ALUAR',9,'LP1912'),
    @('02:58:51','02:58','215_ALUAR',0,'LP1912'),
    @('01:55:51','03:12','215_ALUAR',77,'LP1912'),
    @('03:35:49','04:45','215A_EL PATO',70,'LP1912'),
    @('04:35:25','04:46','215A_EL PATO',11,'LP1912'),
    @('04:48:57','05:34','215B_EL PATO',46,'LP1912'),
    @('05:21:16','05:35','215B_EL PATO',14,'LP1912'),
    @('05:21:16','06:11','215A_EL PATO',50,'LP1912'),
    @('05:52:07','06:12','215A_EL PATO',20,'LP1912'),
    @('06:21:22','06:46','215C_EL PATO',25,'LP1912'),
    @('05:52:07','06:47','215C_EL PATO',55,'LP1912'),
    @('06:59:37','07:11','215A_EL PATO',12,'LP1912'),
    @('05:52:07','07:12','215A_EL PATO',80,'LP1912'),
    @('07:51:34','07:51','215D_EL PATO',0,'LP1912'),
    @('07:51:34','08:23','215B_EL PATO',32,'LP1912'),
    @('07:51:34','09:01','215A_EL PATO',70,'LP1912'),
    @('07:28:14','09:02','215A_EL PATO',94,'LP1912'),
    @('07:51:34','09:42','215C_EL PATO',111,'LP1912')
)
$n2 = $rows2.Count
$arr2 = New-Object 'object[,]' $n2,5
for ($i = 0; $i -lt $n2; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $arr2[$i,$j] = $rows2[$i][$j]
    }
}
$ws.Range($ws.Cells.Item(6,1), $ws.Cells.Item(5+$n2,5)).Value = $arr2

# ---- Sheet3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A2").Value = 'Última actualización: 07:51:34'
$ws.Range("A3").Value = 'Total filas: 19'
$rows3 = @(
    @('00:07:17','00:08','215A_LA PLATA',1,'L6173'),
    @('04:48:57','05:43','215A_LA PLATA',55,'L6173'),
    @('05:21:16','05:44','215A_LA PLATA',23,'L6173'),
    @('04:48:57','06:08','215A_LA PLATA',80,'L6173'),
    @('05:21:16','06:09','215A_LA PLATA',48,'L6173'),
    @('05:52:07','06:13','215A_LA PLATA',21,'L6173'),
    @('04:48:57','06:32','215C_LA PLATA',104,'L6203'),
    @('06:21:22','06:33','215C_LA PLATA',12,'L6203'),
    @('06:59:37','06:59','215B_LP-P MOR-1 Y 57',0,'L6173'),
    @('06:49:33','07:00','215B_LP-P MOR-1 Y 57',11,'L6173'),
    @('07:28:14','07:35','215A_LA PLATA',7,'L6173'),
    @('06:59:37','08:06','215C_LA PLATA',67,'L6203'),
    @('06:49:33','08:07','215C_LA PLATA',78,'L6203'),
    @('07:28:14','08:10','215C_LA PLATA',42,'L6203'),
    @('07:51:34','08:11','215C_LA PLATA',20,'L6203'),
    @('06:49:33','08:33','215A_LA PLATA',104,'L6173'),
    @('07:28:14','08:38','215A_LA PLATA',70,'L6173'),
    @('07:51:34','08:40','215A_LA PLATA',49,'L6173'),
    @('07:51:34','09:09','215D_LA PLATA',78,'L6203')
)
$n3 = $rows3.Count
$arr3 = New-Object 'object[,]' $n3,5
for ($i = 0; $i -lt $n3; $i++) {
    for ($j = 0; $j -lt 5; $j++) {
        $arr3[$i,$j] = $rows3[$i][$j]
    }
}
$ws.Range($ws.Cells.Item(6,1), $ws.Cells.Item(5+$n3,5)).Value = $arr3
